$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J column (k values)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary statistics
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the B14:B17 cells with bold 12pt font and vertical center alignment.
# Apply to one cell first, then propagate via copy/paste-special (formats only)
# to avoid creating redundant intermediate cell styles.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108  # xlCenter
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

[void]$ws.Range("A14:B17").Select()

# Page setup: paper size 9 (A4), portrait orientation
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
